# Refresh the cryptos snapshot table (Sheet1) with the latest scraped
# price/volume figures. Column D holds price text (sometimes genuinely
# numeric-looking, e.g. "563.00") and column E holds the "  +x.xx%  "
# volume-change text. A few rows also got re-ranked, so their Coin name
# (B) and Link (C) cells now hold a different coin entirely.
#
# For price cells that look like plain numbers, Excel's COM layer would
# normally auto-convert the assigned text into a real number (dropping
# the literal formatting, e.g. "563.00" -> 563). To keep these as the
# exact literal text the source data uses, we force the cell to Text
# format before assigning the value, then restore the "Normal" style so
# the cell style stays identical to every other untouched cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '61.985.84'
$ws.Range('E2').Value = '  -0.25%  '
# Row 3
$ws.Range('D3').Value = '2.419.34'
$ws.Range('E3').Value = '  -0.07%  '
# Row 5
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '563.00'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +0.04%  '
# Row 6
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '143.09'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -0.77%  '
# Row 7
$ws.Range('E7').Value = '  -0.07%  '
# Row 8
$ws.Range('E8').Value = '  -0.29%  '
# Row 9
$ws.Range('E9').Value = '  +0.03%  '
# Row 10
$ws.Range('E10').Value = '  -0.83%  '
# Row 11
$ws.Range('E11').Value = '  -3.91%  '
# Row 12
$ws.Range('E12').Value = '  -1.78%  '
# Row 13
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '26.26'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  +0.72%  '
# Row 14
$ws.Range('E14').Value = '  -1.73%  '
# Row 15
$ws.Range('E15').Value = '  -0.39%  '
# Row 16
$ws.Range('D16').Value = '61.890.61'
$ws.Range('E16').Value = '  -0.04%  '
# Row 17
$ws.Range('D17').Value = '2.421.98'
$ws.Range('E17').Value = '  -0.61%  '
# Row 18
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '11.34'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  +1.27%  '
# Row 19
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '323.69'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  -0.26%  '
# Row 20
$ws.Range('E20').Value = '  -1.08%  '
# Row 21
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '6.83'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  +1.08%  '
# Row 22
$ws.Range('E22').Value = '  -0.04%  '
# Row 23
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '66.83'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +2.18%  '
# Row 24
$ws.Range('E24').Value = '  +0.48%  '
# Row 25
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '8.76'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  -3.03%  '
# Row 26
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '552.92'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  -5.80%  '
# Row 27
$ws.Range('B27').Value = 'WrappedeETH'
$ws.Range('C27').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D27').Value = '2.538.02'
$ws.Range('E27').Value = '  +0.44%  '
# Row 28
$ws.Range('B28').Value = 'Binance-PegBSC-USD'
$ws.Range('C28').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '1.00'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  +0.18%  '
# Row 29
$ws.Range('D29').Value = '0.0₃0933'
$ws.Range('E29').Value = '  -0.77%  '
# Row 30
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '8.23'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -0.23%  '
# Row 31
$ws.Range('E31').Value = '  -4.02%  '
# Row 32
$ws.Range('E32').Value = '  -2.05%  '
# Row 33
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '1.87'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  -1.06%  '
# Row 34
$ws.Range('E34').Value = '  -3.24%  '
# Row 35
$ws.Range('E35').Value = '  -0.08%  '
# Row 36
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '4.73'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -0.83%  '
# Row 37
$ws.Range('E37').Value = '  -1.57%  '
# Row 38
$ws.Range('B38').Value = 'RenderToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '5.44'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  -4.76%  '
# Row 39
$ws.Range('B39').Value = 'Monero'
$ws.Range('C39').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '152.42'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  -1.06%  '
# Row 40
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '18.62'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  -0.16%  '
# Row 41
$ws.Range('E41').Value = '  -0.48%  '
# Row 42
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.993'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -0.59%  '
# Row 43
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '147.36'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -1.90%  '
# Row 44
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '2.23'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -5.14%  '
# Row 45
$ws.Range('E45').Value = '  -0.16%  '
# Row 46
$ws.Range('E46').Value = '  -2.02%  '
# Row 47
$ws.Range('E47').Value = '  +0.40%  '
# Row 48
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '19.86'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  -2.33%  '
# Row 49
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.0919'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -0.48%  '
# Row 50
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.0228'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -0.44%  '
# Row 51
$ws.Range('B51').Value = 'BitgetToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/q7gMmMdLb+bitgettoken-bgb'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '1.07'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  +4.31%  '
